$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, "Катамахин",  "Пётр",     "Яковлевич"),
    @(2, "Лазурский",  "Аркадий",  "Виталиевич"),
    @(3, "Вислоушкин", "Семён",    "Данилович"),
    @(4, "Голдякова",  "Ксения",   "Леонидовна"),
    @(5, "Корбуков",   "Степан",   "Борисович"),
    @(6, "Шабуцкий",   "Виталий",  "Юрьевич"),
    @(7, "Морская",    "Лариса",   "Игоревна"),
    @(8, "Лупков",     "Анатолий", "Леонидович"),
    @(9, "Зютина",     "Дарья",    "Леонидовна"),
    @(10, "Зеверов",   "Артём",    "Виталиевич")
)

$rowIndex = 2
foreach ($person in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $person[0]
    $ws.Cells.Item($rowIndex, 2).Value = $person[1]
    $ws.Cells.Item($rowIndex, 3).Value = $person[2]
    $ws.Cells.Item($rowIndex, 4).Value = $person[3]
    $rowIndex++
}
